$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "sql : 3"
$ws.Range("F3").Value = "database : 2"
$ws.Range("F4").Value = "warehouse : 1"
$ws.Range("F5").Value = "design : 2"
$ws.Range("F6").Value = "oracle : 3"
$ws.Range("F7").Value = "documentation : 1"
$ws.Range("F8").Value = "analysis : 1"
